$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2 updates
$ws.Range("G2").Value = 1.65
$ws.Range("H2").Value = 3.8
$ws.Range("I2").Value = 5.25
$ws.Range("J2").Value = 2.3
$ws.Range("O2").Value = 1.4
$ws.Range("P2").Value = 2.75
$ws.Range("X2").Value = 7
$ws.Range("Z2").Value = 12
$ws.Range("AC2").Value = 8
$ws.Range("AD2").Value = 7
$ws.Range("AE2").Value = 21
$ws.Range("AK2").Value = 51
$ws.Range("AL2").Value = 41
$ws.Range("AN2").Value = 3.5
$ws.Range("AO2").Value = 9

# Row 3 updates
$ws.Range("G3").Value = 1.5
$ws.Range("J3").Value = 2.1
$ws.Range("Q3").Value = 2.25
$ws.Range("R3").Value = 1.62
$ws.Range("AE3").Value = 23
$ws.Range("AH3").Value = 13
$ws.Range("AI3").Value = 34
